$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p094r_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p094r_1</id>", 2)
